# Mark Sprint 2 & 3 tasks (rows 6-12) as Done, with start/completed dates
# and developer tracker notes, then update the Sprint Summary roll-up
# numbers for Sprint 2 ("Data & Backend Core") and Sprint 3 ("Today's
# Plan & Activities API").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint & Task Tracking")

# Start/Completed Date columns (L, M) hold plain text dates like "2026-02-12"
# in this workbook (no date number format is defined), so force the cells
# to Text format before assigning the value - this keeps the literal
# "YYYY-MM-DD" string instead of Excel auto-converting it to a date serial.
$dateCols = @("L", "M")
foreach ($col in $dateCols) {
    $ws.Range("${col}6:${col}12").NumberFormat = "@"
}

$rows = @(
    @{ Row = 6;  Start = "2026-02-12"; Completed = "2026-02-12";
       Notes = "TypeORM + PostgreSQL (prod) / better-sqlite3 (dev). 11 entities per §5.1. Indexes on student_id, activity_id, date, syllabus. Versioned migration (idempotent). DatabaseModule with ConfigService." },
    @{ Row = 7;  Start = "2026-02-12"; Completed = "2026-02-12";
       Notes = "10 repositories with BaseRepository error handling: retry 3x for transient errors, 409 for duplicates. All queries parameterized and scoped by student_id. AuthRepository migrated to TypeORM. AuditService persists to DB." },
    @{ Row = 8;  Start = "2026-02-12"; Completed = "2026-02-12";
       Notes = "Business services: StudentService, ActivitiesService, GradingService, AttendanceService, DoubtService. Domain exceptions mapped to HTTP codes. API→Service→Repo layering enforced." },
    @{ Row = 9;  Start = "2026-02-13"; Completed = "2026-02-13";
       Notes = "GET /v1/student/today: student context + task cards (id, type, title, syllabusRef, questionCount, estimatedMinutes, status, score) + progress. Auth+student_id scoped. Empty state: empty tasks, 0%. GET /v1/student/profile: profile + totalActivitiesCompleted. DevSeederService seeds 3 activities+5 questions+10 attendance in dev. Coding standards applied." },
    @{ Row = 10; Start = "2026-02-13"; Completed = "2026-02-13";
       Notes = "GET /v1/student/activities/:type/:id: metadata + questions (correct answer never exposed). PENDING→IN_PROGRESS on first open. POST pause: marks PAUSED, no-op if completed. Types: homework, quiz, test, gap_bridge. 404 if not assigned." },
    @{ Row = 11; Start = "2026-02-13"; Completed = "2026-02-13";
       Notes = "POST /v1/student/activities/:type/:id/respond: {questionId, answer, requestFeedbackLevel?}. Deterministic grading MCQ/TF/FILL_BLANK. Returns isCorrect, score, feedback, feedbackLevel, isComplete, nextQuestionId. Idempotent duplicate submit. Auto-completes on all answered. 409 if completed." },
    @{ Row = 12; Start = "2026-02-13"; Completed = "2026-02-13";
       Notes = "GET /v1/student/results/:type/:id: score, breakdown[{questionId,isCorrect,score}], suggestedNext[{type,title,reason}]. Suggested: <60%→gap_bridge, 60-89%→quiz, 90%+→next topic. 404 if not found." }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Range("K$r").Value = "Done"
    $ws.Range("L$r").Value = $item.Start
    $ws.Range("M$r").Value = $item.Completed
    $ws.Range("P$r").Value = $item.Notes
}

# Sprint Summary roll-up: Sprint 2 planned SP and Sprint 3 capacity/planned/tasks
$ws2 = $wb.Worksheets.Item("Sprint Summary")
$ws2.Range("D3").Value = 3
$ws2.Range("C4").Value = 4
$ws2.Range("D4").Value = 0
$ws2.Range("E4").Value = 0
